# Update the frequency table values in rows 2-5, columns B:X, to reflect
# the "Updated run for publication" results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.00303030303030303,0.00833333333333333,0.975,0.00681818181818182,0.981060606060606,0.000757575757575758,0.00378787878787879,0.0409090909090909,0.00454545454545455,0.999242424242424,0.00454545454545455,0.000757575757575758,0.00227272727272727,0.0053030303030303,0.0121212121212121,0.00681818181818182,0.0106060606060606,0.990151515151515,0.0053030303030303,0.000757575757575758,0.000757575757575758,0.0181818181818182,0.0053030303030303),
    @(0.977272727272727,0.977272727272727,0.00984848484848485,0.0113636363636364,0.0053030303030303,0.00151515151515152,0.85,0.0515151515151515,0.000757575757575758,0,0.0204545454545455,0.986363636363636,0.996212121212121,0.0515151515151515,0,0.00378787878787879,0.000757575757575758,0.000757575757575758,0.00151515151515152,0.968181818181818,0.00227272727272727,0.975757575757576,0.00151515151515152),
    @(0,0.0136363636363636,0.00151515151515152,0.966666666666667,0.0136363636363636,0.00378787878787879,0.0128787878787879,0.00454545454545455,0.99469696969697,0.000757575757575758,0.973484848484849,0,0,0.00151515151515152,0.985606060606061,0.987121212121212,0.975,0.0053030303030303,0.993181818181818,0.00909090909090909,0.00681818181818182,0.00151515151515152,0.0128787878787879),
    @(0.0196969696969697,0.000757575757575758,0.0136363636363636,0.0151515151515152,0,0.993939393939394,0.133333333333333,0.903030303030303,0,0,0,0.0128787878787879,0.00151515151515152,0.941666666666667,0.00151515151515152,0.00151515151515152,0.0136363636363636,0.00378787878787879,0,0.021969696969697,0.990151515151515,0.00454545454545455,0.979545454545454)
)

for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 2
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}
